$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Feb 08 16:44:06 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 16:44:15 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 16:44:26 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 16:44:36 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 16:44:46 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 16:44:56 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 16:45:07 EST 2023"
$ws.Range("B9").Value = "Wed Feb 08 16:45:16 EST 2023"
$ws.Range("B10").Value = "Wed Feb 08 16:45:26 EST 2023"
$ws.Range("B11").Value = "Wed Feb 08 16:45:36 EST 2023"
$ws.Range("B12").Value = "Wed Feb 08 16:45:46 EST 2023"
$ws.Range("B13").Value = "Wed Feb 08 16:45:56 EST 2023"
$ws.Range("B14").Value = "Wed Feb 08 16:46:06 EST 2023"
$ws.Range("B15").Value = "Wed Feb 08 16:46:16 EST 2023"
$ws.Range("B16").Value = "Wed Feb 08 16:46:26 EST 2023"
